$wb = $excel.ActiveWorkbook

# Map of sheet name -> list of (cell, new value) updates for the
# "Elapsed Duration(Hrs)" column (G), reflecting a recalculation of the
# elapsed time since the PCM was created.

$updates = @{
    "R1" = @{ "G2" = "3924:45:23"; "G3" = "64:18:01" }
    "R2" = @{ "G2" = "12106:09:04"; "G3" = "3235:52:33"; "G4" = "474:04:07" }
    "R4" = @{ "G2" = "2951:58:53"; "G3" = "179:11:08" }
    "R5" = @{ "G2" = "425:57:52" }
    "R6" = @{ "G2" = "66:30:10" }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellUpdates = $updates[$sheetName]
    foreach ($cellRef in $cellUpdates.Keys) {
        $ws.Range($cellRef).Value = $cellUpdates[$cellRef]
    }
}
